$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.059.06'
$ws.Range('E2').Value = '  -5.08%  '

$ws.Range('D3').Value = '3.041.24'
$ws.Range('E3').Value = '  -6.68%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').Value = '''548.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -7.25%  '

$ws.Range('D6').Value = '''138.15'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.34%  '

$ws.Range('D7').Value = '''0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.16%  '

$ws.Range('D8').Value = '3.035.66'
$ws.Range('E8').Value = '  -6.59%  '

$ws.Range('D9').Value = '''0.485'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -10.61%  '

$ws.Range('D10').Value = '''6.44'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.97%  '

$ws.Range('D11').Value = '''0.152'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -10.67%  '

$ws.Range('D12').Value = '''0.455'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -9.61%  '

$ws.Range('D13').Value = '''35.60'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -7.06%  '

$ws.Range('D14').Value = '''0.0000218'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -11.26%  '

$ws.Range('D15').Value = '3.527.48'
$ws.Range('E15').Value = '  -6.65%  '

$ws.Range('D16').Value = '64.057.55'
$ws.Range('E16').Value = '  -5.13%  '

$ws.Range('D17').Value = '''0.111'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.57%  '

$ws.Range('D18').Value = '3.038.32'
$ws.Range('E18').Value = '  -6.78%  '

$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '''484.78'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -8.06%  '

$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = '''6.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -7.58%  '

$ws.Range('D21').Value = '''13.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -9.07%  '

$ws.Range('D22').Value = '''0.681'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -9.62%  '

$ws.Range('D23').Value = '''7.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -9.35%  '

$ws.Range('D24').Value = '''77.82'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -8.83%  '

$ws.Range('D25').Value = '''12.34'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -8.39%  '

$ws.Range('D26').Value = '''0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.03%  '

$ws.Range('D27').Value = '''2.70'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -15.87%  '

$ws.Range('D28').Value = '''2.06'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.56%  '

$ws.Range('D29').Value = '''7.68'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.86%  '

$ws.Range('E30').Value = '  -0.21%  '

$ws.Range('D31').Value = '''588.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +14.63%  '

$ws.Range('D32').Value = '''2.63'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.06%  '

$ws.Range('D33').Value = '''25.71'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -11.03%  '

$ws.Range('E34').Value = '  -8.07%  '

$ws.Range('D35').Value = '''5.37'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.35%  '

$ws.Range('D36').Value = '''5.84'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -11.37%  '

$ws.Range('D37').Value = '''51.75'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.48%  '

$ws.Range('D38').Value = '''0.0407'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.54%  '

$ws.Range('D39').Value = '''0.0788'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.38%  '

$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').Value = '''2.79'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.00%  '

$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.118'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.44%  '

$ws.Range('D42').Value = '2.917.01'
$ws.Range('E42').Value = '  -0.43%  '

$ws.Range('D43').Value = '''8.17'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -8.15%  '

$ws.Range('D45').Value = '''0.241'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -9.11%  '

$ws.Range('D46').Value = '''2.06'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.86%  '

$ws.Range('B47').Value = 'PEPE'
$ws.Range('C47').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D47').Value = '0.0₃0529'
$ws.Range('E47').Value = '  -9.80%  '

$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '''24.59'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -7.13%  '

$ws.Range('D49').Value = '''118.54'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.23%  '

$ws.Range('D50').Value = '''0.108'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.29%  '

$ws.Range('D51').Value = '''2.03'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -12.05%  '
